$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before existing row 285, shifting the old data (rows
# 285-299) down to rows 288-302.
$ws.Rows.Item(285).Insert()
$ws.Rows.Item(285).Insert()
$ws.Rows.Item(285).Insert()

function Set-DataRow {
    param(
        [int]$Row,
        $Fecha,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [double]$PrecioKg,
        [double]$KgUnidad
    )

    $ws.Cells.Item($Row, 1).Value = 5
    $ws.Cells.Item($Row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($Row, 3).Value = "Maule"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($Row, 5).Value = 7
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value = 100103
    $ws.Cells.Item($Row, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($Row, 9).Value = 100103001
    $ws.Cells.Item($Row, 10).Value = "Cereza"
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = "Provincia de Curicó"
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

# New row 285
Set-DataRow 285 "2023-01-05" "Bing" "Primera" 220 3500 4000 3727 "$/bandeja 10 kilos" 373 10

# New row 286
Set-DataRow 286 "2023-01-05" "Lapins" "Primera" 170 3000 4000 3706 "$/bandeja 10 kilos" 371 10

# New row 287
Set-DataRow 287 "2023-01-05" "Lapins" "Segunda" 80 2500 2500 2500 "$/bandeja 10 kilos" 250 10

# Existing row (old 285) shifted to 288 also changes its Fecha value.
$ws.Cells.Item(288, 4).Value = "2023-01-05"
$ws.Cells.Item(288, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
